$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$s.Shapes.Item(4).Delete()
